$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values - B2 and D2 are cleared, C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -9.093568270444603
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -9.4901985549390577

# Row 3 values updated
$ws.Range("B3").Value = -12.183051192106124
$ws.Range("C3").Value = -3.9144188612369137
$ws.Range("D3").Value = -16.548441345687586
$ws.Range("E3").Value = 19.710913143782037

# Update the selection to reflect the new active range
$ws.Range("B1:E3").Select()
